$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 3 with the new incoming mail ---
$logs = $wb.Worksheets.Item("Logs")
$logs.Range("A3").Value = "Probleem met inloggen"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Range("D3").Value = "IT / Technisch probleem"
$logs.Range("F3").Value = "2025-06-22 17:08:10"
$logs.Range("G3").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too
$catFc = $logs.Range("D2").FormatConditions.Item(1)
$catFc.ModifyAppliesToRange($logs.Range("D2:D3"))

$answeredFc = $logs.Range("G2").FormatConditions.Item(1)
$answeredFc.ModifyAppliesToRange($logs.Range("G2:G3"))

# --- Dashboard sheet: add the category count row for the new category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A3").Value = "IT / Technisch probleem"
$dash.Range("B3").Value = 1

# --- Chart: extend the category/value series references to include row 3 ---
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
